$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 0.09138
$ws.Range("M2").Value = 0.6746743333333334
$ws.Range("N2").Value = 2.024023
$ws.Range("O2").Value = 0.07069047851636343
$ws.Range("P2").Value = 0.07069047851636343
$ws.Range("Q2").Value = 0.02055058019333333
$ws.Range("R2").Value = 0.18495522174
$ws.Range("S2").Value = 0.07069047851636343
$ws.Range("T2").Value = 0.07069047851636343

# Row 3
$ws.Range("H3").Value = 0.09138
$ws.Range("O3").Value = 0.4692497642600617
$ws.Range("P3").Value = 0.4692497642600616
$ws.Range("S3").Value = 0.4692497642600617
$ws.Range("T3").Value = 0.4692497642600616

# Row 4
$ws.Range("H4").Value = 0.09138
$ws.Range("O4").Value = 0.460059757223575
$ws.Range("P4").Value = 0.460059757223575
$ws.Range("S4").Value = 0.460059757223575
$ws.Range("T4").Value = 0.460059757223575
